$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text format on Price cells whose new values would otherwise be
# auto-detected as numbers by Excel, so they stay stored as literal text
# (matching the source formatting, e.g. "583.42" as a string, not a number).
$ws.Range("D5:D7").NumberFormat = "@"
$ws.Range("D9:D13").NumberFormat = "@"
$ws.Range("D18:D33").NumberFormat = "@"
$ws.Range("D35:D36").NumberFormat = "@"
$ws.Range("D39:D40").NumberFormat = "@"
$ws.Range("D42:D44").NumberFormat = "@"
$ws.Range("D46:D48").NumberFormat = "@"
$ws.Range("D50:D51").NumberFormat = "@"

# Apply the updated cell values
$ws.Range("D2").Value = "66.906.81"
$ws.Range("E2").Value = "  -3.32%  "
$ws.Range("D3").Value = "3.591.47"
$ws.Range("E3").Value = "  -3.82%  "
$ws.Range("E4").Value = "  +0.65%  "
$ws.Range("D5").Value = "583.42"
$ws.Range("E5").Value = "  -3.58%  "
$ws.Range("D6").Value = "182.87"
$ws.Range("E6").Value = "  -2.02%  "
$ws.Range("D7").Value = "0.603"
$ws.Range("E7").Value = "  -5.32%  "
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("D9").Value = "0.668"
$ws.Range("E9").Value = "  -7.84%  "
$ws.Range("D10").Value = "0.142"
$ws.Range("E10").Value = "  -13.81%  "
$ws.Range("D11").Value = "53.05"
$ws.Range("E11").Value = "  -7.05%  "
$ws.Range("D12").Value = "0.0000246"
$ws.Range("E12").Value = "  -17.05%  "
$ws.Range("D13").Value = "9.83"
$ws.Range("E13").Value = "  -6.82%  "
$ws.Range("D14").Value = "4.177.40"
$ws.Range("E14").Value = "  -3.04%  "
$ws.Range("D15").Value = "3.595.16"
$ws.Range("E15").Value = "  -3.16%  "
$ws.Range("E16").Value = "  -0.77%  "
$ws.Range("D17").Value = "66.657.93"
$ws.Range("E17").Value = "  -3.20%  "
$ws.Range("D18").Value = "18.24"
$ws.Range("E18").Value = "  -6.79%  "
$ws.Range("D19").Value = "12.11"
$ws.Range("E19").Value = "  -6.80%  "
$ws.Range("D20").Value = "1.05"
$ws.Range("E20").Value = "  -6.92%  "
$ws.Range("D21").Value = "391.10"
$ws.Range("E21").Value = "  -5.64%  "
$ws.Range("D22").Value = "4.28"
$ws.Range("E22").Value = "  -8.50%  "
$ws.Range("D23").Value = "84.57"
$ws.Range("E23").Value = "  -5.90%  "
$ws.Range("D24").Value = "2.81"
$ws.Range("E24").Value = "  -8.02%  "
$ws.Range("D25").Value = "6.04"
$ws.Range("E25").Value = "  -0.37%  "
$ws.Range("D26").Value = "12.09"
$ws.Range("E26").Value = "  -6.90%  "
$ws.Range("D27").Value = "10.17"
$ws.Range("E27").Value = "  -8.22%  "
$ws.Range("D28").Value = "3.57"
$ws.Range("E28").Value = "  -9.85%  "
$ws.Range("D29").Value = "8.83"
$ws.Range("E29").Value = "  -7.95%  "
$ws.Range("D30").Value = "30.84"
$ws.Range("E30").Value = "  -7.11%  "
$ws.Range("D31").Value = "6.69"
$ws.Range("E31").Value = "  -9.39%  "
$ws.Range("D32").Value = "65.74"
$ws.Range("E32").Value = "  +1.13%  "
$ws.Range("D33").Value = "11.75"
$ws.Range("E33").Value = "  -6.78%  "
$ws.Range("E34").Value = "  -6.54%  "
$ws.Range("D35").Value = "577.20"
$ws.Range("E35").Value = "  -5.61%  "
$ws.Range("D36").Value = "41.16"
$ws.Range("E36").Value = "  -7.22%  "
$ws.Range("E37").Value = "  -0.13%  "
$ws.Range("E38").Value = "  +0.45%  "
$ws.Range("D39").Value = "0.369"
$ws.Range("E39").Value = "  -8.93%  "
$ws.Range("D40").Value = "0.131"
$ws.Range("E40").Value = "  -4.63%  "
$ws.Range("D41").Value = "0.0₃0715"
$ws.Range("E41").Value = "  -20.92%  "
$ws.Range("D42").Value = "2.74"
$ws.Range("E42").Value = "  -11.04%  "
$ws.Range("D43").Value = "0.0405"
$ws.Range("E43").Value = "  -8.74%  "
$ws.Range("B44").Value = "Stellar"
$ws.Range("C44").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D44").Value = "0.130"
$ws.Range("E44").Value = "  -4.28%  "
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Value = "2.662.32"
$ws.Range("E45").Value = "  -4.89%  "
$ws.Range("B46").Value = "ApeXProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D46").Value = "3.04"
$ws.Range("E46").Value = "  -3.21%  "
$ws.Range("B47").Value = "Monero"
$ws.Range("C47").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D47").Value = "139.45"
$ws.Range("E47").Value = "  -1.73%  "
$ws.Range("B48").Value = "Fetch.AI"
$ws.Range("C48").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D48").Value = "2.33"
$ws.Range("E48").Value = "  -16.30%  "
$ws.Range("E49").Value = "  -7.80%  "
$ws.Range("D50").Value = "8.33"
$ws.Range("E50").Value = "  -10.91%  "
$ws.Range("D51").Value = "2.53"
$ws.Range("E51").Value = "  -9.65%  "
